$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Sector of Institution" values copied in from Arrest.xls for
# Reporting Location = On Campus, Offense = Arrest - Illegal Weapon Possessions.
$sectors = @(
    "Public, 4-year or above",
    "Private nonprofit, 4-year or above",
    "Private for-profit, 4-year or above",
    "Public, 2-year",
    "Private nonprofit, 2-year",
    "Private for-profit, 2-year",
    "Public, less-than 2-year",
    "Private nonprofit, less-than 2-year",
    "Private for-profit, less-than 2-year"
)

# Write them into A2:A10, using a leading apostrophe so Excel stores the
# cells as explicit text (quote-prefixed), matching the source workbook.
for ($i = 0; $i -lt $sectors.Length; $i++) {
    $row = $i + 2
    $ws.Range("A$row").Value = "'" + $sectors[$i]
}

# Column sizing to fit the newly populated data.
$ws.Columns.Item(1).ColumnWidth = 28
$ws.Columns.Item(2).ColumnWidth = 15.65

# Restore the active selection to B7, as recorded in the source file.
$ws.Range("B7").Select() | Out-Null
